$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entries: 25.10.18 (row 16) and 26.10.18 (row 17)
# Row 16 - "25. 10.18", 13:00 -> 14:30, remark "-Issues with WMI"
$ws.Range("A16").Value = "25. 10.18"
$ws.Range("B16").Value = 0.54166666666666663
$ws.Range("C16").Value = 0.60416666666666663
$ws.Range("E16").Value = "-Issues with WMI"

# Row 17 - "26. 10.18", 18:30 -> 20:00, remark "-Issues with WMI"
$ws.Range("A17").Value = "26. 10.18"
$ws.Range("B17").Value = 0.77083333333333337
$ws.Range("C17").Value = 0.83333333333333337
$ws.Range("E17").Value = "-Issues with WMI"

# Move the active selection to C18, matching the saved cursor position
$ws.Range("C18").Select()
